$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) + 2) Relocate the bold "Play Astro Cat Free ..." paragraph that
#    used to sit right before the closing italic summary paragraph:
#    cut it from the bottom of the document and paste it in right
#    after the very first (Heading1) paragraph. This removes it from
#    the bottom (satisfying the deletion in the diff) and reuses its
#    exact run/paragraph formatting (no pStyle, leading empty run,
#    bold run) for the new "Meta description" paragraph at the top.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$oldHeadingPara = $d.Paragraphs.Item($count - 1)
$oldHeadingPara.Range.Cut() | Out-Null

$p1 = $d.Paragraphs.First
$pastePos = $p1.Range.End
$pasteRange = $d.Range($pastePos, $pastePos)
$pasteRange.Paste() | Out-Null

# The pasted paragraph now reads "Play Astro Cat Free - Innovative
# gameplay with 1296 ways to win" in bold - turn it into the new
# "Meta description" lead-in, keeping the bold run, then append the
# (non-bold) description sentence as a second run.
$newPara = $d.Paragraphs.Item(2)
$boldStart = $newPara.Range.Start
$boldEnd = $newPara.Range.End

$boldText = "Meta description"
$boldTextRange = $d.Range($boldStart, $boldEnd)
$boldTextRange.Text = $boldText

$restText = ": Read our review of Astro Cat, a well-designed online slot game with an innovative gameplay structure, exciting bonuses, and 1296 ways to win. Play for free now."
$restPos = $boldStart + $boldText.Length
$restInsert = $d.Range($restPos, $restPos)
$restInsert.InsertAfter($restText) | Out-Null
$restRange = $d.Range($restPos, $restPos + $restText.Length)
$restRange.Font.Bold = $false

# ------------------------------------------------------------------
# 3) Replace the text of the remaining (now last) italic paragraph
#    with the new feature-image design-prompt copy. Assign directly
#    to Range.Text (not Find.Execute's replacement argument) so the
#    straight quotes in the new copy are kept verbatim instead of
#    being auto-converted to curly/smart quotes.
# ------------------------------------------------------------------
$newSummary = 'Create a feature image for "Astro Cat" Design Prompt: Create a cartoon-style feature image for "Astro Cat" that features a happy Maya warrior with glasses. Requirements: - The image must be colorful and eye-catching - The happy Maya warrior with glasses should be the main focus of the image - The background should have a space or cosmic theme - The image should incorporate elements from the game, such as the Astro Cat wild symbol or the fireworks bonus symbol Suggestions: - The happy Maya warrior could be holding a large Astro Cat symbol, or sitting on a pile of coins won from the game - The background could have a starry sky or galaxies to enhance the cosmic theme - The image could feature other characters or symbols from the game, such as the koi carp or Chinese fan, to add more detail and depth to the image.'

$lastCount = $d.Paragraphs.Count
$summaryPara = $d.Paragraphs.Item($lastCount)
$summaryRange = $summaryPara.Range
$summaryRange.MoveEnd(1, -1) | Out-Null
$summaryRange.Text = $newSummary

Write-Output "done"
